# Apply the PlayerPerformance_4685.xlsx edit:
#  1. Remove the stray empty B3 cell on the "ODI Batting" sheet.
#  2. Add a new "ODI Batting Extra" worksheet (after "ODI Bowling") with
#     match/batting-position/boundary/percentage/MOM data.

$wb = $excel.ActiveWorkbook

# --- 1. Clear the orphan empty inline string cell at ODI Batting!B3 ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B3").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet at the end of the tab strip ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row: reuse the same bold/centered/bordered header style already
# used by the other sheets' row 1 (copy format from an existing header).
$wb.Worksheets.Item("ODI Bowling").Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# Row 2
$extra.Range("A2").Value = "'4023"
$extra.Range("B2").Value = 8
$extra.Range("C2").Value = "'1"
$extra.Range("D2").Value = "'0"
$extra.Range("E2").Value = "'3.81%"
$extra.Range("F2").Value = "NO"

# Row 3 (NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are blank for this entry, but
# the cells themselves are still present in the source data, not removed)
$extra.Range("A3").Value = "'4026"
$extra.Range("B3").Value = 9
$extra.Range("C3").Value = "'"
$extra.Range("D3").Value = "'"
$extra.Range("E3").Value = "'"
$extra.Range("F3").Value = "NO"

# The leading apostrophes above force these numeric-looking values (match
# codes, counts, percentages) to be stored as text instead of being
# auto-coerced to numbers, same as the source data. That quote-prefix entry
# mode tags the cells with a "quote prefix" style in Excel; reset the style
# back to Normal so only the cell content/type changes, not its formatting.
$extra.Range("A2:A3").Style = "Normal"
$extra.Range("C2:E3").Style = "Normal"
